$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 969.04
$ws.Range("I19").Value = 946.6667
$ws.Range("J19").Value = 1002.6
$ws.Range("K19").Value = 946.6667
$ws.Range("L19").Value = 1002.6
$ws.Range("M19").Value = -771.6667
$ws.Range("N19").Value = -1352.6

$ws.Range("H33").Value = 4844.1816
$ws.Range("I33").Value = 6562.8125
$ws.Range("J33").Value = 261.16666
$ws.Range("K33").Value = 6562.8125
$ws.Range("L33").Value = 261.16666
$ws.Range("M33").Value = -6333.8125
$ws.Range("N33").Value = -719.16666

$ws.Range("L63").ClearContents()
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("N63").Value = 0

$ws.Range("L66").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("N66").Value = 0

$ws.Range("H80").Value = 1193
$ws.Range("I80").Value = 608.6
$ws.Range("J80").Value = 1436.5
$ws.Range("K80").Value = 1825.8
$ws.Range("L80").Value = 4309.5
$ws.Range("M80").Value = -827.8000000000002
$ws.Range("N80").Value = -6305.5

$ws.Range("H83").Value = 1193
$ws.Range("I83").Value = 608.6
$ws.Range("J83").Value = 1436.5
$ws.Range("K83").Value = 5477.400000000001
$ws.Range("L83").Value = 12928.5
$ws.Range("M83").Value = -485.4000000000005
$ws.Range("N83").Value = -22912.5

$ws.Range("H132").Value = 3415.9
$ws.Range("I132").Value = 2945.2144
$ws.Range("J132").Value = 10005.5
$ws.Range("K132").Value = 8835.643199999999
$ws.Range("L132").Value = 30016.5
$ws.Range("M132").Value = -6305.643199999999
$ws.Range("N132").Value = -35076.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3014
$ws.Range("I2").Value = 2172.9375
$ws.Range("J2").Value = 8060.375
$ws.Range("K2").Value = 2172.9375
$ws.Range("L2").Value = 8060.375
$ws.Range("M2").Value = -2059.9375
$ws.Range("N2").Value = -8286.375

$ws.Range("H32").Value = 2707.7
$ws.Range("I32").Value = 1789.6511
$ws.Range("J32").Value = 8347.143
$ws.Range("K32").Value = 1789.6511
$ws.Range("L32").Value = 8347.143
$ws.Range("M32").Value = -1502.6511
$ws.Range("N32").Value = -8921.143

$ws.Range("H61").Value = 3973.762
$ws.Range("I61").Value = 1672.9412
$ws.Range("J61").Value = 13752.25
$ws.Range("K61").Value = 1672.9412
$ws.Range("L61").Value = 13752.25
$ws.Range("M61").Value = -1460.9412
$ws.Range("N61").Value = -14176.25

$ws.Range("L114").ClearContents()
$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("N114").Value = 0

$ws.Range("H116").Value = 3014
$ws.Range("I116").Value = 2172.9375
$ws.Range("J116").Value = 8060.375
$ws.Range("K116").Value = 2172.9375
$ws.Range("L116").Value = 8060.375
$ws.Range("M116").Value = 121.0625
$ws.Range("N116").Value = -12648.375

$ws.Range("H122").Value = 4500
$ws.Range("I122").Value = 4500
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 13500
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -11050

$ws.Range("H129").Value = 56155.2
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 56155.2
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 56155.2
$ws.Range("N129").Value = -66155.2

$ws.Range("H132").Value = 6792.4
$ws.Range("I132").Value = 5591.912
$ws.Range("J132").Value = 12056.077
$ws.Range("K132").Value = 16775.736
$ws.Range("L132").Value = 36168.231
$ws.Range("M132").Value = -14245.736
$ws.Range("N132").Value = -41228.231

$ws.Range("H136").Value = 3973.762
$ws.Range("I136").Value = 1672.9412
$ws.Range("J136").Value = 13752.25
$ws.Range("K136").Value = 5018.8236
$ws.Range("L136").Value = 41256.75
$ws.Range("M136").Value = -2468.8236
$ws.Range("N136").Value = -46356.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3014
$ws.Range("I3").Value = 2172.9375
$ws.Range("J3").Value = 8060.375
$ws.Range("K3").Value = 2172.9375
$ws.Range("L3").Value = 8060.375
$ws.Range("M3").Value = -2058.9375
$ws.Range("N3").Value = -8288.375

$ws.Range("H94").Value = 513.26086
$ws.Range("I94").Value = 513.26086
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 513.26086
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -62.26085999999998

$ws.Range("H105").Value = 3041.3076
$ws.Range("I105").Value = 3209.8948
$ws.Range("J105").Value = 2583.7144
$ws.Range("K105").Value = 3209.8948
$ws.Range("L105").Value = 2583.7144
$ws.Range("M105").Value = -1462.8948
$ws.Range("N105").Value = -6077.7144

$ws.Range("H134").Value = 7895.3794
$ws.Range("I134").Value = 5162.8
$ws.Range("J134").Value = 24974
$ws.Range("K134").Value = 15488.4
$ws.Range("L134").Value = 74922
$ws.Range("M134").Value = -12953.4
$ws.Range("N134").Value = -79992

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2987.375
$ws.Range("I31").Value = 7449.5
$ws.Range("J31").Value = 1500
$ws.Range("K31").Value = 7449.5
$ws.Range("L31").Value = 1500
$ws.Range("M31").Value = -7154.5
$ws.Range("N31").Value = -2090

$ws.Range("H34").Value = 2987.375
$ws.Range("I34").Value = 7449.5
$ws.Range("J34").Value = 1500
$ws.Range("K34").Value = 7449.5
$ws.Range("L34").Value = 1500
$ws.Range("M34").Value = -7247.5
$ws.Range("N34").Value = -1904

$ws.Range("H56").Value = 40449.5
$ws.Range("I56").Value = 50899
$ws.Range("J56").Value = 30000
$ws.Range("K56").Value = 50899
$ws.Range("L56").Value = 30000
$ws.Range("M56").Value = -50054
$ws.Range("N56").Value = -31690

$ws.Range("H94").Value = 828.61536
$ws.Range("I94").Value = 566.8333
$ws.Range("J94").Value = 1053
$ws.Range("K94").Value = 566.8333
$ws.Range("L94").Value = 1053
$ws.Range("M94").Value = -115.8333
$ws.Range("N94").Value = -1955

$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("N122").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1873.375
$ws.Range("I5").Value = 1152.7273
$ws.Range("J5").Value = 2754.1667
$ws.Range("K5").Value = 3458.1819
$ws.Range("L5").Value = 8262.500100000001
$ws.Range("M5").Value = -3346.1819
$ws.Range("N5").Value = -8486.500100000001

$ws.Range("H11").Value = 600.3333
$ws.Range("I11").Value = 600.3333
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1800.9999
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -1660.9999

$ws.Range("H40").Value = 126
$ws.Range("I40").Value = 172.25
$ws.Range("J40").Value = 79.75
$ws.Range("K40").Value = 689
$ws.Range("L40").Value = 319
$ws.Range("M40").Value = -620
$ws.Range("N40").Value = -457

$ws.Range("H80").Value = 23288.4
$ws.Range("I80").Value = 32600.4
$ws.Range("J80").Value = 13976.4
$ws.Range("K80").Value = 97801.20000000001
$ws.Range("L80").Value = 41929.2
$ws.Range("M80").Value = -96865.20000000001
$ws.Range("N80").Value = -43801.2

$ws.Range("H83").Value = 23288.4
$ws.Range("I83").Value = 32600.4
$ws.Range("J83").Value = 13976.4
$ws.Range("K83").Value = 293403.6
$ws.Range("L83").Value = 125787.6
$ws.Range("M83").Value = -288723.6
$ws.Range("N83").Value = -135147.6

$ws.Range("H101").Value = 12000
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 12000
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 36000
$ws.Range("N101").Value = -40868

$ws.Range("H135").Value = 1873.375
$ws.Range("I135").Value = 1152.7273
$ws.Range("J135").Value = 2754.1667
$ws.Range("K135").Value = 10374.5457
$ws.Range("L135").Value = 24787.5003
$ws.Range("M135").Value = -7839.545700000001
$ws.Range("N135").Value = -29857.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 4174166.2
$ws.Range("I11").Value = 6253750
$ws.Range("J11").Value = 14999
$ws.Range("K11").Value = 6253750
$ws.Range("L11").Value = 14999
$ws.Range("M11").Value = -6253611
$ws.Range("N11").Value = -15277

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3475827
$ws.Range("I136").Value = 4632841.5
$ws.Range("J136").Value = 4783.375
$ws.Range("K136").Value = 13898524.5
$ws.Range("L136").Value = 14350.125
$ws.Range("M136").Value = -13895974.5
$ws.Range("N136").Value = -19450.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5887.7144
$ws.Range("I122").Value = 6257.8887
$ws.Range("J122").Value = 3666.6667
$ws.Range("K122").Value = 18773.6661
$ws.Range("L122").Value = 11000.0001
$ws.Range("M122").Value = -16323.6661
$ws.Range("N122").Value = -15900.0001

$ws.Range("H136").Value = 3803.6316
$ws.Range("I136").Value = 3718.577
$ws.Range("J136").Value = 3987.9167
$ws.Range("K136").Value = 11155.731
$ws.Range("L136").Value = 11963.7501
$ws.Range("M136").Value = -8605.731
$ws.Range("N136").Value = -17063.7501
